$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# ---------------------------------------------------------------------------
# 1. Insert a new "version list" sheet right after "Export as TSV" and
#    before "assay_category list".
# ---------------------------------------------------------------------------
$versionSheet = $wb.Worksheets.Add($null, $ws)
$versionSheet.Name = "version list"
$versionSheet.Range("A1").NumberFormat = "@"
$versionSheet.Range("A1").Value = "1"

# ---------------------------------------------------------------------------
# 2. On the "Export as TSV" sheet: remove the existing comments (they do not
#    travel with a column insert), then insert two new columns in front of
#    the data (shifting everything - headers, data and validations - two
#    columns to the right).
# ---------------------------------------------------------------------------
while ($ws.Comments.Count -gt 0) {
    $ws.Comments.Item(1).Delete()
}

$ws.Range("A:B").EntireColumn.Insert(-4161)

# ---------------------------------------------------------------------------
# 3. Populate the two new header cells and copy the header formatting
#    (bold, centered, wrapped) from the neighbouring header cell.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "version"
$ws.Range("B1").Value = "description"

$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Add the data validation for the new "version" column.
# ---------------------------------------------------------------------------
$rngVersion = $ws.Range("A2:A1048576")
$rngVersion.Validation.Add(3, 1, 1, "='version list'!`$A`$1:`$A`$1") | Out-Null
$rngVersion.Validation.ErrorTitle = "Value must come from list"
$rngVersion.Validation.ErrorMessage = "Value must be one of: 1."
$rngVersion.Validation.ShowInput = $true
$rngVersion.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 5. Re-create every header comment at its (shifted) location.
# ---------------------------------------------------------------------------
$comments = @{
    "A1" = "Version of the schema to use when validating this metadata.";
    "B1" = "Free-text description of this assay.";
    "C1" = "HuBMAP Display ID of the donor of the assayed tissue.";
    "D1" = "HuBMAP Display ID of the assayed tissue.";
    "E1" = "Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.";
    "F1" = "DOI for protocols.io referring to the protocol for this assay.";
    "G1" = "Name of the person responsible for executing the assay.";
    "H1" = "Email address for the operator.";
    "I1" = "Name of the principal investigator responsible for the data.";
    "J1" = "Email address for the principal investigator.";
    "K1" = "Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.";
    "L1" = "The specific type of assay being executed.";
    "M1" = "Analytes are the target molecules being measured with the assay.";
    "N1" = "Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein.";
    "O1" = "An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.";
    "P1" = "Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.";
    "Q1" = "The width of a pixel.";
    "R1" = "The unit of measurement of the width of a pixel.";
    "S1" = "The height of a pixel";
    "T1" = "The unit of measurement of the height of a pixel.";
    "U1" = "Optional if assay does not have multiple z-levels. Note that this is resolution within a given sample: z-pitch (resolution_z_value) is the increment distance between image slices ie. the microscope stage is moved up or down in increments to capture images of several focal planes. The best one will be used & the rest discarded. The thickness of the sample itself is sample metadata. Leave blank if not applicable.";
    "V1" = "The unit of incremental distance between image slices.";
    "W1" = "Chemical stains (dyes) applied to histology samples to highlight important features of the tissue as well as to enhance the tissue contrast.";
    "X1" = "DOI for protocols.io referring to the protocol for preparing tissue sections for the assay.";
    "Y1" = "DOI for protocols.io for the overall process.";
    "Z1" = "Relative path to file with ORCID IDs for contributors for this dataset.";
    "AA1" = "Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions."
}

$order = @("A1","B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1","M1","N1","O1","P1","Q1","R1","S1","T1","U1","V1","W1","X1","Y1","Z1","AA1")
foreach ($ref in $order) {
    $ws.Range($ref).AddComment($comments[$ref]) | Out-Null
}

$ws.Range("A1").Select() | Out-Null
